# Update the "取得日時" (acquired timestamp) column on the "ランサーズ" sheet
# for all data rows (2-26) from "2026-02-04 01:52:02" to "2026-02-04 02:26:18".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "2026-02-04 01:52:02"
$newValue = "2026-02-04 02:26:18"

for ($row = 2; $row -le 26; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
